$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.345.36'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.846.85'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9979'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '240.24'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.6265'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.9989'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.07590'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -1.24%  '
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -1.49%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07736'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.6787'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -0.25%  '
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -2.17%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '82.94'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -0.95%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '6.122'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -0.54%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '29.382.10'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '227.71'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -0.67%  '
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -1.06%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.9987'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.468'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.9988'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '158.55'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +0.82%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '8.426'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +0.54%  '
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.432'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +8.94%  '
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -0.69%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.05604'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -2.00%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.066'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +0.27%  '
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -1.13%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.6958'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -1.84%  '
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -0.25%  '
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.227.26'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +0.27%  '
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -2.09%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.354'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -1.50%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.8963'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -1.48%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.9986'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '101.40'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -0.44%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '65.47'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -1.20%  '
$ws.Range('B44').Value = 'BabyDogeCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.00000000122'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +1.07%  '
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '7.213'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +0.96%  '
$ws.Range('B46').Value = 'TheSandbox'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.3985'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -1.02%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.027'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.688'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.1140'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +1.32%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.05696'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.4621'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -0.16%  '
